$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 318: update date and price columns ---
$ws.Range("D318").Value = 44628
$ws.Range("K318").Value = 450
$ws.Range("L318").Value = 500
$ws.Range("M318").Value = 475
$ws.Range("P318").Value = 475

# --- Row 319: update date and price columns ---
$ws.Range("D319").Value = 44628
$ws.Range("K319").Value = 350
$ws.Range("L319").Value = 400
$ws.Range("M319").Value = 375
$ws.Range("P319").Value = 375

# --- Row 320: update date and quality ---
$ws.Range("D320").Value = 44552
$ws.Range("I320").Value = "Segunda"

# --- Row 321: update date, quality, volume, price columns ---
$ws.Range("D321").Value = 44552
$ws.Range("I321").Value = "Tercera"
$ws.Range("J321").Value = 1200
$ws.Range("K321").Value = 300
$ws.Range("L321").Value = 350
$ws.Range("M321").Value = 325
$ws.Range("P321").Value = 325

# --- Row 322: update quality and price columns ---
$ws.Range("I322").Value = "Primera"
$ws.Range("K322").Value = 400
$ws.Range("L322").Value = 450
$ws.Range("M322").Value = 425
$ws.Range("P322").Value = 425

# --- Row 323: update date, volume, price columns ---
$ws.Range("D323").Value = 44544
$ws.Range("J323").Value = 1300
$ws.Range("K323").Value = 350
$ws.Range("L323").Value = 400
$ws.Range("M323").Value = 375
$ws.Range("P323").Value = 375

# --- Row 324: update date and price columns ---
$ws.Range("D324").Value = 44544
$ws.Range("L324").Value = 350
$ws.Range("M324").Value = 325
$ws.Range("P324").Value = 325

# --- New Row 325 (copy of former row 323 content: 44160 / Segunda) ---
$ws.Range("A325").Value = 1
$ws.Range("B325").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C325").Value = "Arica y Parinacota"
$ws.Range("D325").Value = 44160
$ws.Range("E325").Value = 15
$ws.Range("F325").Value = 100112023
$ws.Range("G325").Value = "Brócoli"
$ws.Range("H325").Value = "Sin especificar"
$ws.Range("I325").Value = "Segunda"
$ws.Range("J325").Value = 1200
$ws.Range("K325").Value = 400
$ws.Range("L325").Value = 450
$ws.Range("M325").Value = 425
$ws.Range("N325").Value = "$/unidad"
$ws.Range("O325").Value = "Región de Arica y Parinacota"
$ws.Range("P325").Value = 425
$ws.Range("Q325").Value = 1
$ws.Range("R325").Value = "Hortaliza"

# --- New Row 326 (copy of former row 324 content: 44160 / Tercera) ---
$ws.Range("A326").Value = 1
$ws.Range("B326").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C326").Value = "Arica y Parinacota"
$ws.Range("D326").Value = 44160
$ws.Range("E326").Value = 15
$ws.Range("F326").Value = 100112023
$ws.Range("G326").Value = "Brócoli"
$ws.Range("H326").Value = "Sin especificar"
$ws.Range("I326").Value = "Tercera"
$ws.Range("J326").Value = 1200
$ws.Range("K326").Value = 300
$ws.Range("L326").Value = 400
$ws.Range("M326").Value = 350
$ws.Range("N326").Value = "$/unidad"
$ws.Range("O326").Value = "Región de Arica y Parinacota"
$ws.Range("P326").Value = 350
$ws.Range("Q326").Value = 1
$ws.Range("R326").Value = "Hortaliza"

# Apply the same date display format used elsewhere (style s="2") to the new D cells
$ws.Range("D325").NumberFormat = $ws.Range("D324").NumberFormat
$ws.Range("D326").NumberFormat = $ws.Range("D324").NumberFormat
